# Fixed Pathing issue when printing
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "52" "56"
Replace-Text "Megrendelő: Példa Péter" "Megrendelő: Kelemen Kabátban"
Replace-Text "Cím: 9700 Szombathely Nincs Ilyen út 69" "Cím: 9702 Szomszédfalu Valamien utca 3"
Replace-Text "Elérhetőség: telefon  06301234567" "Elérhetőség: telefon  062056473829"
Replace-Text "Megjegyzés:Valamien megjegyzés" "Megjegyzés:"
Replace-Text "Megnevezés: Fûnyíró" "Megnevezés: Traktor"
Replace-Text "Típus: Husqwarna" "Típus: EpikusFunyirok"
Replace-Text "Modell: CW23" "Modell: Sututu3"
Replace-Text "Hibajelenség: Nem indul el" "Hibajelenség: Nem forog a kerék"
Replace-Text "Tartozékok: Kerék, Kesztyû" "Tartozékok: Az egész traktor"
Replace-Text "Szerviz diagnózis: El van törve" "Szerviz diagnózis: Ki kell engedni a féket"
